$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.848.50"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.304.08"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'306.90"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").Value = "'96.55"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").Value = "'35.43"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'18.46"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D15").Value = "2.662.63"
$ws.Range("D16").Value = "2.318.30"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "42.784.58"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'13.06"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'67.37"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'4.01"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'25.32"
$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  +16.56%  "
$ws.Range("D30").Value = "'166.36"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D32").Value = "'33.20"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'4.76"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D36").Value = "'17.77"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "'0.0693"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "2.012.06"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'18.27"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("D46").Value = "'10.06"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "'2.07"
$ws.Range("E47").Value = "  -7.34%  "
$ws.Range("D48").Value = "'2.81"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  +10.68%  "
$ws.Range("D50").Value = "'53.86"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "2.528.43"
$ws.Range("E51").Value = "  -0.07%  "
